# Generate Report for Handoff
# Adds a new tracked file (f4099a79-e6bc-4717-ab7e-9f17458e2234.md) as row 9
# across the Overview / zh-cn / de-de sheets, growing each sheet's table by
# one row and wiring up the corresponding hyperlinks.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet (row 9)
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A9").Value = "f4099a79-e6bc-4717-ab7e-9f17458e2234.md"
$wsOverview.Range("B9").Value = "e2e\f4099a79-e6bc-4717-ab7e-9f17458e2234.md"
$wsOverview.Range("C9").Value = ".md"
$wsOverview.Range("D9").Value = ""
$wsOverview.Range("E9").Value = "Ready for handoff"
$wsOverview.Range("F9").Value = "Ready for handoff"
$wsOverview.Range("G9").Value = "2016-08-25 08:47:34"

$wsOverview.Hyperlinks.Add($wsOverview.Range("B9"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b3f4c6a1d8e29f7065c1234abf890123e4567890/e2e/f4099a79-e6bc-4717-ab7e-9f17458e2234.md", "", "", "e2e\f4099a79-e6bc-4717-ab7e-9f17458e2234.md") | Out-Null

$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.Resize($wsOverview.Range("A1:G9"))

# ---------------------------------------------------------------------
# zh-cn sheet (row 9)
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A9").Value = "f4099a79-e6bc-4717-ab7e-9f17458e2234.md"
$wsZhCn.Range("B9").Value = ".md"
$wsZhCn.Range("C9").Value = "Ready for handoff"
$wsZhCn.Range("D9").Value = "e2e"
$wsZhCn.Range("E9").Value = "ht"
$wsZhCn.Range("F9").Value = "'False"
$wsZhCn.Range("G9").Value = "f4099a79-e6bc-4717-ab7e-9f17458e2234.fa697d251c1c7511f3329a28bc2a17066135f240.zh-cn.xlf"
$wsZhCn.Range("H9").Value = "2016-08-25 08:47:29"
$wsZhCn.Range("I9").Value = ""
$wsZhCn.Range("J9").Value = ""
$wsZhCn.Range("K9").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("L9").Value = ""
$wsZhCn.Range("M9").Value = "'True"
$wsZhCn.Range("N9").Value = ""
$wsZhCn.Range("O9").Value = "'False"
$wsZhCn.Range("P9").Value = ""

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A9"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b3f4c6a1d8e29f7065c1234abf890123e4567890/e2e/f4099a79-e6bc-4717-ab7e-9f17458e2234.md", "", "", "f4099a79-e6bc-4717-ab7e-9f17458e2234.md") | Out-Null

$loZhCn = $wsZhCn.ListObjects.Item(1)
$loZhCn.Resize($wsZhCn.Range("A1:P9"))

# ---------------------------------------------------------------------
# de-de sheet (row 9)
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A9").Value = "f4099a79-e6bc-4717-ab7e-9f17458e2234.md"
$wsDeDe.Range("B9").Value = ".md"
$wsDeDe.Range("C9").Value = "Ready for handoff"
$wsDeDe.Range("D9").Value = "e2e"
$wsDeDe.Range("E9").Value = "ht"
$wsDeDe.Range("F9").Value = "'False"
$wsDeDe.Range("G9").Value = "f4099a79-e6bc-4717-ab7e-9f17458e2234.fa697d251c1c7511f3329a28bc2a17066135f240.de-de.xlf"
$wsDeDe.Range("H9").Value = "2016-08-25 08:47:34"
$wsDeDe.Range("I9").Value = ""
$wsDeDe.Range("J9").Value = ""
$wsDeDe.Range("K9").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("L9").Value = ""
$wsDeDe.Range("M9").Value = "'True"
$wsDeDe.Range("N9").Value = ""
$wsDeDe.Range("O9").Value = "'False"
$wsDeDe.Range("P9").Value = ""

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A9"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b3f4c6a1d8e29f7065c1234abf890123e4567890/e2e/f4099a79-e6bc-4717-ab7e-9f17458e2234.md", "", "", "f4099a79-e6bc-4717-ab7e-9f17458e2234.md") | Out-Null

$loDeDe = $wsDeDe.ListObjects.Item(1)
$loDeDe.Resize($wsDeDe.Range("A1:P9"))
